# Update NATMI LR-pair TPM-derived metrics with recomputed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = [double]"0.3333333333333333"
$ws.Range("G2").Value = [double]"0.002991666666666667"
$ws.Range("H2").Value = [double]"0.008975"
$ws.Range("I2").Value = [double]"0.0003566413595017623"
$ws.Range("J2").Value = [double]"0.0003566413595017623"
$ws.Range("M2").Value = [double]"3.112844666666666"
$ws.Range("N2").Value = [double]"9.338534"
$ws.Range("O2").Value = [double]"0.0962303687181678"
$ws.Range("P2").Value = [double]"0.09623036871816783"
$ws.Range("Q2").Value = [double]"0.009312593627777777"
$ws.Range("R2").Value = [double]"0.08381334265"
$ws.Range("S2").Value = [double]"3.431972952500323e-05"
$ws.Range("T2").Value = [double]"3.431972952500323e-05"

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = [double]"0.3333333333333333"
$ws.Range("G3").Value = [double]"0.002991666666666667"
$ws.Range("H3").Value = [double]"0.008975"
$ws.Range("I3").Value = [double]"0.0003566413595017623"
$ws.Range("J3").Value = [double]"0.0003566413595017623"
$ws.Range("O3").Value = [double]"0.5562258596073433"
$ws.Range("P3").Value = [double]"0.5562258596073434"
$ws.Range("Q3").Value = [double]"0.05382817778611111"
$ws.Range("R3").Value = [double]"0.484453600075"
$ws.Range("S3").Value = [double]"0.0001983731467603993"
$ws.Range("T3").Value = [double]"0.0001983731467603993"

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = [double]"0.3333333333333333"
$ws.Range("G4").Value = [double]"0.002991666666666667"
$ws.Range("H4").Value = [double]"0.008975"
$ws.Range("I4").Value = [double]"0.0003566413595017623"
$ws.Range("J4").Value = [double]"0.0003566413595017623"
$ws.Range("M4").Value = [double]"11.24229066666667"
$ws.Range("N4").Value = [double]"33.726872"
$ws.Range("O4").Value = [double]"0.3475437716744887"
$ws.Range("P4").Value = [double]"0.3475437716744888"
$ws.Range("Q4").Value = [double]"0.03363318624444445"
$ws.Range("R4").Value = [double]"0.3026986762"
$ws.Range("S4").Value = [double]"0.0001239484832163597"
$ws.Range("T4").Value = [double]"0.0001239484832163598"

# Row 5
$ws.Range("I5").Value = [double]"0.9971069332391614"
$ws.Range("J5").Value = [double]"0.9971069332391616"
$ws.Range("M5").Value = [double]"3.112844666666666"
$ws.Range("N5").Value = [double]"9.338534"
$ws.Range("O5").Value = [double]"0.0962303687181678"
$ws.Range("P5").Value = [double]"0.09623036871816783"
$ws.Range("Q5").Value = [double]"26.03638480312088"
$ws.Range("R5").Value = [double]"234.327463228088"
$ws.Range("S5").Value = [double]"0.09595196783704603"
$ws.Range("T5").Value = [double]"0.09595196783704607"

# Row 6
$ws.Range("I6").Value = [double]"0.9971069332391614"
$ws.Range("J6").Value = [double]"0.9971069332391616"
$ws.Range("O6").Value = [double]"0.5562258596073433"
$ws.Range("P6").Value = [double]"0.5562258596073434"
$ws.Range("S6").Value = [double]"0.5546166610613944"
$ws.Range("T6").Value = [double]"0.5546166610613946"

# Row 7
$ws.Range("I7").Value = [double]"0.9971069332391614"
$ws.Range("J7").Value = [double]"0.9971069332391616"
$ws.Range("M7").Value = [double]"11.24229066666667"
$ws.Range("N7").Value = [double]"33.726872"
$ws.Range("O7").Value = [double]"0.3475437716744887"
$ws.Range("P7").Value = [double]"0.3475437716744888"
$ws.Range("Q7").Value = [double]"94.03251276887822"
$ws.Range("R7").Value = [double]"846.2926149199039"
$ws.Range("S7").Value = [double]"0.3465383043407208"
$ws.Range("T7").Value = [double]"0.3465383043407209"

# Row 8
$ws.Range("G8").Value = [double]"0.02127666666666667"
$ws.Range("H8").Value = [double]"0.06383"
$ws.Range("I8").Value = [double]"0.002536425401336767"
$ws.Range("J8").Value = [double]"0.002536425401336767"
$ws.Range("M8").Value = [double]"3.112844666666666"
$ws.Range("N8").Value = [double]"9.338534"
$ws.Range("O8").Value = [double]"0.0962303687181678"
$ws.Range("P8").Value = [double]"0.09623036871816783"
$ws.Range("Q8").Value = [double]"0.06623095835777777"
$ws.Range("R8").Value = [double]"0.5960786252199999"
$ws.Range("S8").Value = [double]"0.0002440811515967639"
$ws.Range("T8").Value = [double]"0.0002440811515967639"

# Row 9
$ws.Range("G9").Value = [double]"0.02127666666666667"
$ws.Range("H9").Value = [double]"0.06383"
$ws.Range("I9").Value = [double]"0.002536425401336767"
$ws.Range("J9").Value = [double]"0.002536425401336767"
$ws.Range("O9").Value = [double]"0.5562258596073433"
$ws.Range("P9").Value = [double]"0.5562258596073434"
$ws.Range("Q9").Value = [double]"0.3828248009011111"
$ws.Range("R9").Value = [double]"3.44542320811"
$ws.Range("S9").Value = [double]"0.001410825399188444"
$ws.Range("T9").Value = [double]"0.001410825399188444"

# Row 10
$ws.Range("G10").Value = [double]"0.02127666666666667"
$ws.Range("H10").Value = [double]"0.06383"
$ws.Range("I10").Value = [double]"0.002536425401336767"
$ws.Range("J10").Value = [double]"0.002536425401336767"
$ws.Range("M10").Value = [double]"11.24229066666667"
$ws.Range("N10").Value = [double]"33.726872"
$ws.Range("O10").Value = [double]"0.3475437716744887"
$ws.Range("P10").Value = [double]"0.3475437716744888"
$ws.Range("Q10").Value = [double]"0.2391984710844444"
$ws.Range("R10").Value = [double]"2.15278623976"
$ws.Range("S10").Value = [double]"0.000881518850551559"
$ws.Range("T10").Value = [double]"0.0008815188505515591"

Write-Host "Updated 100 cells"